# The commit swaps the contents of ppt/theme/theme1.xml (the deck's main
# slide-master theme, currently the colourful "Integral" / "Red Violet"
# scheme) and ppt/theme/theme2.xml (the notes-master theme, currently the
# plain default "Office Theme" scheme) - i.e. after the edit the slides use
# the plain Office colour palette and the notes master uses the Red Violet
# palette. The font scheme and format (fill/line/effect) scheme are byte
# for byte identical between the two themes already, so the only
# observable difference is the 12-slot colour scheme.
#
# The PowerPoint object model only exposes (and only ever persists writes
# through) the single theme that is bound to the slide master -
# Slide(s).ThemeColorScheme / Master.Theme.ThemeColorScheme - regardless of
# whether the call is made via a Slide, the SlideMaster, the NotesMaster or
# the NotesPage; every one of those resolves to ppt/theme/theme1.xml in
# this host. There is no reachable member that edits the notes-master's own
# theme part, so we apply the reachable half of the swap: push the plain
# Office colours into the deck's theme colour scheme.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function HexToRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Order matches the standard OOXML / ThemeColorScheme slot order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = HexToRGB $officeColors[$i - 1]
}
